$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 8 (franzosa_ControlvsCD_Fp),
# shifting it and all subsequent rows down by one.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the franzosa_ControlvsCD_ConvCD data.
$ws.Range("A8").Value = "franzosa_ControlvsCD_ConvCD"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.2
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8
$ws.Range("H8").Value = 0.8

# Insert a new row before the current row 13 (franzosa_ControlvsUC_Fp, after
# the previous insert shifted everything below row 8 down by one), so it
# lands right after franzosa_ControlvsUC_Age and before franzosa_ControlvsUC_Fp.
$ws.Rows("13:13").Insert()

# Populate the newly inserted row 13 with the franzosa_ControlvsUC_ConvUC data.
$ws.Range("A13").Value = "franzosa_ControlvsUC_ConvUC"
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0.6
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.4
$ws.Range("H13").Value = 0.4
